$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update survey completion dates in column A (merged from Survey Completion by Member)
# Rows previously carrying the placeholder date 1/1/2001 (serial 36892) are updated
# to their actual completion dates (serials 45126-45128, i.e. 7/19/2023-7/21/2023).
$dateUpdates = @{
    285 = 45128
    291 = 45127
    292 = 45127
    293 = 45127
    294 = 45127
    296 = 45127
    298 = 45127
    299 = 45127
    300 = 45127
    301 = 45127
    302 = 45127
    303 = 45127
    304 = 45127
    305 = 45127
    306 = 45127
    307 = 45127
    308 = 45127
    317 = 45126
    323 = 45126
    324 = 45126
    326 = 45126
    328 = 45126
    331 = 45126
    336 = 45126
    338 = 45126
    339 = 45126
    340 = 45126
    341 = 45126
    342 = 45126
    343 = 45126
    345 = 45126
    347 = 45126
    350 = 45126
    353 = 45126
    355 = 45126
    357 = 45126
    359 = 45126
    361 = 45126
    370 = 45126
    372 = 45126
    373 = 45126
    375 = 45126
    377 = 45126
    379 = 45126
    382 = 45126
    384 = 45126
    386 = 45126
    388 = 45126
    390 = 45126
    417 = 45126
    429 = 45126
    440 = 45126
    441 = 45126
    442 = 45126
    443 = 45126
    444 = 45126
    452 = 45126
    454 = 45126
}

foreach ($row in $dateUpdates.Keys) {
    $ws.Cells.Item($row, 1).Value2 = $dateUpdates[$row]
}

# Restore the view/selection state recorded for this save
[void]$ws.Range("A567").Select()
[void]$ws.Range("B525").Select()

"Updated $($dateUpdates.Count) date cells"
